$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 636.03845
$ws.Range("J17").Value = 658.96
$ws.Range("L17").Value = 1976.88
$ws.Range("N17").Value = -2312.88
$ws.Range("H32").Value = 691.7857
$ws.Range("I32").Value = 570.75
$ws.Range("K32").Value = 570.75
$ws.Range("M32").Value = -244.75
$ws.Range("H43").Value = 588998.4399999999
$ws.Range("I43").Value = 1996.6666
$ws.Range("J43").Value = 1029249.75
$ws.Range("K43").Value = 1996.6666
$ws.Range("L43").Value = 1029249.75
$ws.Range("M43").Value = -1927.6666
$ws.Range("N43").Value = -1029387.75
$ws.Range("H88").Value = 72703540
$ws.Range("I88").Value = 277778270
$ws.Range("J88").Value = 11181109
$ws.Range("K88").Value = 277778270
$ws.Range("L88").Value = 11181109
$ws.Range("M88").Value = -277777864
$ws.Range("N88").Value = -11181921
$ws.Range("H91").Value = 72703540
$ws.Range("I91").Value = 277778270
$ws.Range("J91").Value = 11181109
$ws.Range("K91").Value = 277778270
$ws.Range("L91").Value = 11181109
$ws.Range("M91").Value = -277776866
$ws.Range("N91").Value = -11183917
$ws.Range("H137").Value = 3447.5557
$ws.Range("I137").Value = 5111
$ws.Range("J137").Value = 2807.7693
$ws.Range("K137").Value = 15333
$ws.Range("L137").Value = 8423.3079
$ws.Range("M137").Value = -12783
$ws.Range("N137").Value = -13523.3079

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 4999.6665
$ws.Range("J11").Value = 6500
$ws.Range("L11").Value = 6500
$ws.Range("N11").Value = -6788
$ws.Range("H32").Value = 2234484.2
$ws.Range("I32").Value = 2502152.5
$ws.Range("K32").Value = 2502152.5
$ws.Range("M32").Value = -2501865.5
$ws.Range("H74").Value = 39630.742
$ws.Range("I74").Value = 49239.76
$ws.Range("K74").Value = 49239.76
$ws.Range("M74").Value = -48365.76
$ws.Range("H77").Value = 39630.742
$ws.Range("I77").Value = 49239.76
$ws.Range("K77").Value = 246198.8
$ws.Range("M77").Value = -241830.8
$ws.Range("H132").Value = 5599.447
$ws.Range("I132").Value = 4102.3237
$ws.Range("J132").Value = 9515
$ws.Range("K132").Value = 12306.9711
$ws.Range("L132").Value = 28545
$ws.Range("M132").Value = -9776.971099999999
$ws.Range("N132").Value = -33605

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 20835586
$ws.Range("I20").Value = 23811848
$ws.Range("J20").Value = 1748
$ws.Range("K20").Value = 23811848
$ws.Range("L20").Value = 1748
$ws.Range("M20").Value = -23811601
$ws.Range("N20").Value = -2242
$ws.Range("H107").Value = 26787034
$ws.Range("I107").Value = 34092170
$ws.Range("K107").Value = 34092170
$ws.Range("M107").Value = -34090250

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 1500
$ws.Range("J13").Value = 1500
$ws.Range("L13").Value = 1500
$ws.Range("N13").Value = -1778
$ws.Range("H31").Value = 9360.467000000001
$ws.Range("I31").Value = 4217.75
$ws.Range("K31").Value = 4217.75
$ws.Range("M31").Value = -3922.75
$ws.Range("H34").Value = 9360.467000000001
$ws.Range("I34").Value = 4217.75
$ws.Range("K34").Value = 4217.75
$ws.Range("M34").Value = -4015.75
$ws.Range("H58").Value = 7497.6665
$ws.Range("J58").Value = 8403.629999999999
$ws.Range("L58").Value = 8403.629999999999
$ws.Range("N58").Value = -8809.629999999999
$ws.Range("H70").Value = 65000
$ws.Range("J70").Value = 65000
$ws.Range("L70").Value = 65000
$ws.Range("N70").Value = -65630
$ws.Range("H73").Value = 65000
$ws.Range("J73").Value = 65000
$ws.Range("L73").Value = 65000
$ws.Range("N73").Value = -67184
$ws.Range("H86").Value = 30313684
$ws.Range("I86").Value = 10113502
$ws.Range("K86").Value = 10113502
$ws.Range("M86").Value = -10112379
$ws.Range("H89").Value = 30313684
$ws.Range("I89").Value = 10113502
$ws.Range("K89").Value = 50567510
$ws.Range("M89").Value = -50561894
$ws.Range("H99").Value = 5724.3335
$ws.Range("I99").Value = 4000
$ws.Range("K99").Value = 4000
$ws.Range("M99").Value = -2502
$ws.Range("H100").Value = 39500
$ws.Range("J100").Value = 39500
$ws.Range("L100").Value = 39500
$ws.Range("N100").Value = -41664
$ws.Range("H107").Value = 1243.7727
$ws.Range("I107").Value = 827.3461
$ws.Range("K107").Value = 827.3461
$ws.Range("M107").Value = 1092.6539
$ws.Range("H122").Value = 2610.2104
$ws.Range("I122").Value = 2387.8823
$ws.Range("K122").Value = 7163.646900000001
$ws.Range("M122").Value = -4713.646900000001
$ws.Range("H126").Value = 5724.3335
$ws.Range("I126").Value = 4000
$ws.Range("K126").Value = 12000
$ws.Range("M126").Value = -9530
$ws.Range("H132").Value = 6055.4863
$ws.Range("I132").Value = 3968.3914
$ws.Range("K132").Value = 11905.1742
$ws.Range("M132").Value = -9375.174199999999
$ws.Range("H134").Value = 10132.206
$ws.Range("I134").Value = 12282
$ws.Range("K134").Value = 36846
$ws.Range("M134").Value = -34311
$ws.Range("H136").Value = 7497.6665
$ws.Range("J136").Value = 8403.629999999999
$ws.Range("L136").Value = 25210.89
$ws.Range("N136").Value = -30310.89

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2974.5356
$ws.Range("J113").Value = 4069.6843
$ws.Range("L113").Value = 12209.0529
$ws.Range("N113").Value = -16549.0529
$ws.Range("H131").Value = 1730.1875
$ws.Range("I131").Value = 948.6667
$ws.Range("K131").Value = 2846.0001
$ws.Range("M131").Value = 2193.9999
$ws.Range("H140").Value = 252943.5
$ws.Range("I140").Value = 402009.6
$ws.Range("K140").Value = 1206028.8
$ws.Range("M140").Value = -1200848.8

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 250
$ws.Range("I7").Value = 250
$ws.Range("K7").Value = 250
$ws.Range("M7").Value = -138
$ws.Range("H8").Value = 250
$ws.Range("I8").Value = 250
$ws.Range("K8").Value = 250
$ws.Range("M8").Value = -111

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 83335720
$ws.Range("I68").Value = 142858620
$ws.Range("K68").Value = 142858620
$ws.Range("M68").Value = -142857871
$ws.Range("H71").Value = 83335720
$ws.Range("I71").Value = 142858620
$ws.Range("K71").Value = 714293100
$ws.Range("M71").Value = -714289356
$ws.Range("H82").Value = 4999.6665
$ws.Range("I82").Value = 4999.5
$ws.Range("J82").Value = 5000
$ws.Range("K82").Value = 4999.5
$ws.Range("L82").Value = 5000
$ws.Range("M82").Value = -4638.5
$ws.Range("N82").Value = -5722
$ws.Range("H85").Value = 4999.6665
$ws.Range("I85").Value = 4999.5
$ws.Range("J85").Value = 5000
$ws.Range("K85").Value = 4999.5
$ws.Range("L85").Value = 5000
$ws.Range("M85").Value = -3751.5
$ws.Range("N85").Value = -7496
$ws.Range("H100").Value = 3960.6
$ws.Range("J100").Value = 4001.3333
$ws.Range("L100").Value = 4001.3333
$ws.Range("N100").Value = -5083.3333
$ws.Range("H132").Value = 6388.7617
$ws.Range("J132").Value = 11114.125
$ws.Range("L132").Value = 33342.375
$ws.Range("N132").Value = -38402.375
$ws.Range("H136").Value = 10080.125
$ws.Range("I136").Value = 4810.5
$ws.Range("J136").Value = 11836.667
$ws.Range("K136").Value = 14431.5
$ws.Range("L136").Value = 35510.001
$ws.Range("M136").Value = -11881.5
$ws.Range("N136").Value = -40610.001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8690.1
$ws.Range("I62").Value = 7900.5
$ws.Range("K62").Value = 7900.5
$ws.Range("M62").Value = -7276.5
$ws.Range("H65").Value = 8690.1
$ws.Range("I65").Value = 7900.5
$ws.Range("K65").Value = 39502.5
$ws.Range("M65").Value = -36382.5
